$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New theta_out values (direct/reverse simulation re-run), asymptoting to 80
$newValues = @(
    79.96902787125583,
    79.96985792333928,
    79.97064429775531,
    79.97141624200937,
    79.97215815361176,
    79.97286133728616,
    79.97356333979303,
    79.97422380008364,
    79.97486107913733,
    79.97547705574212,
    79.97607188859901,
    79.97665026258868,
    79.97719034981064,
    79.97772506806807,
    79.97823300163626,
    79.97872402641761,
    79.97919584714261,
    79.97965436686627,
    79.98008546660564,
    79.98049628724532,
    79.98089756232307,
    79.98128237541196,
    79.98165453962814,
    79.98201172670397,
    79.9823584350162,
    79.98268151508753,
    79.9830008690334,
    79.98329964829878,
    79.98359355295378,
    79.98387010762542,
    79.98414058576236,
    79.98439351636232,
    79.98464504912135,
    79.984885048865,
    79.98510538191158,
    79.98533151785071,
    79.98553811143243,
    79.98574546087599,
    79.98593215702415,
    79.98611806235354,
    79.98630087445098,
    79.9864734702071,
    79.98663465141932,
    79.98679089527701,
    79.9869465057443,
    79.98709072922364,
    79.98723241290011,
    79.98736719737354,
    79.9874928465792,
    79.98762116102495,
    79.9877386526538,
    79.9878590371339,
    79.98796596096106,
    79.98807512998822,
    79.98817742146245,
    79.9882779856873,
    79.98837278939909,
    79.98846177585176,
    79.98855020760837,
    79.9886385757041,
    79.98871545718862,
    79.98879637466773,
    79.9888755387708,
    79.9889497298991,
    79.98901819872016,
    79.9890826070821,
    79.9891545340509,
    79.98921656311173,
    79.98928022962362,
    79.98934028252648,
    79.9893947215607,
    79.98945361934801,
    79.98950613560901,
    79.98956074593261,
    79.98960398922722,
    79.98965569962104,
    79.98970500735992,
    79.98974941597301,
    79.98979710181891,
    79.98983800647162,
    79.98988357602131,
    79.9899222979161,
    79.98995930821471,
    79.98999625328959,
    79.99003203655293,
    79.9900704922034,
    79.99010913112652,
    79.9901428181705,
    79.99017165436132,
    79.99020421298677,
    79.99023665454362,
    79.99026671261042,
    79.99029678889376,
    79.99032661433951,
    79.99035726715091,
    79.99038363294643,
    79.99040916597531,
    79.99043195801318,
    79.99045827487133,
    79.99048662531818
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $newValues[$i]
}

